# Updated cryptos list on Mon Oct 30 17:47:05 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    # Forces the written value to stay a text string (matching the original
    # inlineStr cells) even when its content looks like a number, without
    # leaving a visible style index on the cell afterwards.
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "34.410.43"
$ws.Range("E2").Value = "  -0.57%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.800.48"
$ws.Range("E3").Value = "  +0.01%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.53%  "

# Row 5 - BNB
Set-TextValue "D5" "227.62"
$ws.Range("E5").Value = "  +0.29%  "

# Row 6 - XRP
Set-TextValue "D6" "0.579"
$ws.Range("E6").Value = "  +3.54%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.50%  "

# Row 8 - Solana
Set-TextValue "D8" "34.86"
$ws.Range("E8").Value = "  +5.60%  "

# Row 9 - Cardano
Set-TextValue "D9" "0.298"
$ws.Range("E9").Value = "  +0.46%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -0.63%  "

# Row 11 - TRON
Set-TextValue "D11" "0.0952"
$ws.Range("E11").Value = "  +0.28%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "2.060.28"
$ws.Range("E12").Value = "  +0.28%  "

# Row 13 - Chainlink
Set-TextValue "D13" "11.15"

# Row 14 - WrappedEther
$ws.Range("D14").Value = "1.790.68"
$ws.Range("E14").Value = "  +0.04%  "

# Row 15 - Polygon
$ws.Range("E15").Value = "  +0.23%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "34.370.56"
$ws.Range("E16").Value = "  -0.41%  "

# Row 17 - Polkadot
$ws.Range("E17").Value = "  +1.17%  "

# Row 18 - Litecoin
Set-TextValue "D18" "69.22"
$ws.Range("E18").Value = "  +0.34%  "

# Row 19 - ShibaInu
$ws.Range("D19").Value = "0.0₃0797"
$ws.Range("E19").Value = "  -0.65%  "

# Row 20 - BitcoinCash
Set-TextValue "D20" "244.77"
$ws.Range("E20").Value = "  -1.56%  "

# Row 21 - Avalanche
Set-TextValue "D21" "11.50"
$ws.Range("E21").Value = "  +0.93%  "

# Row 23 - Uniswap
Set-TextValue "D23" "4.14"
$ws.Range("E23").Value = "  -0.82%  "

# Row 24 - Monero
Set-TextValue "D24" "170.93"
$ws.Range("E24").Value = "  +3.99%  "

# Row 25 - Toncoin
$ws.Range("E25").Value = "  +2.11%  "

# Row 26 - Cosmos
Set-TextValue "D26" "7.54"
$ws.Range("E26").Value = "  +3.55%  "

# Row 27 - EthereumClassic
Set-TextValue "D27" "16.71"
$ws.Range("E27").Value = "  +0.82%  "

# Row 28 - Stellar
$ws.Range("E28").Value = "  +1.65%  "

# Row 29 - BinanceUSD
$ws.Range("E29").Value = "  +0.33%  "

# Row 30 - InternetComputer(DFINITY)
$ws.Range("E30").Value = "  +1.15%  "

# Row 31/32 - Hedera and PancakeSwap swap places
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue "D31" "1.25"
$ws.Range("E31").Value = "  +0.80%  "

$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D32" "0.0527"
$ws.Range("E32").Value = "  +0.99%  "

# Row 33 - Filecoin
$ws.Range("E33").Value = "  -0.30%  "

# Row 34 - LidoDAOToken
$ws.Range("E34").Value = "  -0.06%  "

# Row 35 - Maker
$ws.Range("D35").Value = "1.399.78"

# Row 36 - RenderToken
Set-TextValue "D36" "2.54"
$ws.Range("E36").Value = "  -1.77%  "

# Row 37 - ImmutableX
$ws.Range("E37").Value = "  +0.60%  "

# Row 38 - TrustWalletToken
$ws.Range("E38").Value = "  -0.37%  "

# Row 39 - VeChain
$ws.Range("E39").Value = "  -2.04%  "

# Row 40 - Aave
$ws.Range("E40").Value = "  -2.98%  "

# Row 41 - MXToken
Set-TextValue "D41" "2.83"
$ws.Range("E41").Value = "  +3.10%  "

# Row 42 - ARBITRUM
Set-TextValue "D42" "0.948"
$ws.Range("E42").Value = "  +0.75%  "

# Row 43 - HuobiToken
$ws.Range("E43").Value = "  +0.43%  "

# Row 44 - InjectiveProtocol
Set-TextValue "D44" "13.61"
$ws.Range("E44").Value = "  +0.93%  "

# Row 45 - WEMIXToken
Set-TextValue "D45" "1.10"
$ws.Range("E45").Value = "  +2.83%  "

# Row 46 - Kaspa
Set-TextValue "D46" "0.0512"
$ws.Range("E46").Value = "  -2.00%  "

# Row 47 - FraxShare
$ws.Range("E47").Value = "  -1.14%  "

# Row 48 - RocketPoolETH
$ws.Range("D48").Value = "1.962.16"
$ws.Range("E48").Value = "  +0.51%  "

# Row 49 - Quant
Set-TextValue "D49" "104.39"
$ws.Range("E49").Value = "  -1.40%  "

# Row 50 - PaxDollar
$ws.Range("E50").Value = "  +0.48%  "

# Row 51 - BabyDogeCoin
$ws.Range("D51").Value = "0.0₆0129"
$ws.Range("E51").Value = "  +0.09%  "
